# Apply product-name lookups to the "Facturados" sheet (invoice/remision rows)
# and update the resulting price/subtotal values. Also reset the active tab
# back to the first sheet ("Pendientes").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Facturados")

# Row 2: "Pedrito" / "Clavos 100gr"
$ws.Range("G2").Value = "Pedrito"
$ws.Range("O2").Value = "Clavos 100gr"
$ws.Range("R2").Value = 2000
$ws.Range("S2").Value = 360000

# Row 3: "Juanito" / "Clavos 200gr"
$ws.Range("G3").Value = "Juanito"
$ws.Range("O3").Value = "Clavos 200gr"
$ws.Range("R3").Value = 3000
$ws.Range("S3").Value = 612000

# Make "Pendientes" (first sheet) the active tab again.
$wb.Worksheets.Item("Pendientes").Activate()
